$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-1650996087652503"
$ws1.Range("B2").Value = "go_stims-16509960876124754.csv"
$ws1.Range("B3").Value = "GNG_stims-16509960876365027.csv"
$ws1.Range("B4").Value = "go_stims-16509960876365027.csv"
$ws1.Range("B5").Value = "GNG_stims-1650996087652503.csv"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-1650996089884481"
$ws2.Range("B2").Value = "OB-1650996088428465.csv"
$ws2.Range("B3").Value = "ZB-match_3-1650996087956463.csv"
$ws2.Range("B4").Value = "ZB-match_9-16509960880924666.csv"
$ws2.Range("B5").Value = "OB-16509960889164999.csv"
$ws2.Range("B6").Value = "TB-16509960890924957.csv"
$ws2.Range("B7").Value = "TB-16509960898604972.csv"
$ws2.Range("B8").Value = "ZB-match_9-16509960879004633.csv"
$ws2.Range("B9").Value = "TB-16509960893885045.csv"
$ws2.Range("B10").Value = "OB-16509960889404688.csv"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-1650996089884481"

$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16509960899405038"
$ws4.Range("B2").Value = "MM_stims-1650996089908499.csv"
$ws4.Range("B3").Value = "ZM_stims-1650996089884481.csv"
$ws4.Range("B4").Value = "MM_stims-16509960899244711.csv"
$ws4.Range("B5").Value = "ZM_stims-1650996089908499.csv"
$ws4.Range("B6").Value = "MM_stims-16509960899405038.csv"
$ws4.Range("B7").Value = "ZM_stims-16509960899244711.csv"

$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16509960900045"
$ws5.Range("B2").Value = "SAT_stims-16509960899405038.csv"
$ws5.Range("B3").Value = "SAT_stims-16509960899564643.csv"
$ws5.Range("B4").Value = "vSAT_stims-1650996089988464.csv"
$ws5.Range("B5").Value = "vSAT_stims-16509960899724967.csv"
